# "selesai vr3 eps 2" - fill in the results for the 2nd epsilon-reduction run (row 7)
# and move the active selection to F10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F7").Value = 72
$ws.Range("G7").Value = 227
$ws.Range("H7").Value = "59m 12.97s"

$ws.Range("F10").Select()
